$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.82"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.87"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.368"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05635"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.433"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.355"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8159"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9306"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07504"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03247"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03088"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09326"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.563"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001591"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04726"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005783"

$ws.Range("E18").Value = "17OneONE"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006379"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005059"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001035"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.747"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.157"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003002"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03938"

$ws.Range("B41").Value = "BKEXToken"

$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1066"

$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003022"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002921"

$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005575"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7804"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1794"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
